# Fix formatting on the fastq "purpose" column: fullRNASEQ -> fullRNASeq
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
